$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'60.149.56"
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").Formula = "'2.417.02"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Formula = "'552.46"
$ws.Range("E5").Value = "  -0.50%  "

$ws.Range("D6").Formula = "'137.17"
$ws.Range("E6").Value = "  -1.54%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Formula = "'0.587"
$ws.Range("E8").Value = "  +1.63%  "

$ws.Range("E9").Value = "  -1.31%  "

$ws.Range("D10").Formula = "'5.66"
$ws.Range("E10").Value = "  -1.90%  "

$ws.Range("E11").Value = "  -0.27%  "

$ws.Range("E12").Value = "  -1.58%  "

$ws.Range("D13").Formula = "'24.88"
$ws.Range("E13").Value = "  -0.51%  "

$ws.Range("D14").Formula = "'2.849.51"
$ws.Range("E14").Value = "  -0.48%  "

$ws.Range("D15").Formula = "'60.034.65"
$ws.Range("E15").Value = "  -0.29%  "

$ws.Range("E16").Value = "  -1.48%  "

$ws.Range("D17").Formula = "'2.432.85"
$ws.Range("E17").Value = "  +0.05%  "

$ws.Range("E18").Value = "  -0.63%  "

$ws.Range("E19").Value = "  +1.26%  "

$ws.Range("D20").Formula = "'328.00"
$ws.Range("E20").Value = "  -1.78%  "

$ws.Range("E21").Value = "  -0.39%  "

$ws.Range("D22").Formula = "'1.00"
$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").Formula = "'65.47"
$ws.Range("E23").Value = "  +0.37%  "

$ws.Range("E24").Value = "  +3.00%  "

$ws.Range("E25").Value = "  +0.50%  "

$ws.Range("D26").Formula = "'1.01"
$ws.Range("E26").Value = "  +0.62%  "

$ws.Range("E27").Value = "  +2.11%  "

$ws.Range("E28").Value = "  -2.07%  "

$ws.Range("E29").Value = "  -2.28%  "

$ws.Range("D30").Formula = "'170.03"
$ws.Range("E30").Value = "  +0.38%  "

$ws.Range("D31").Formula = "'6.09"
$ws.Range("E31").Value = "  -3.96%  "

$ws.Range("E32").Value = "  +2.79%  "

$ws.Range("E33").Value = "  -4.20%  "

$ws.Range("D34").Formula = "'18.56"
$ws.Range("E34").Value = "  -1.15%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").Formula = "'1.31"
$ws.Range("E36").Value = "  +0.42%  "

$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("E38").Value = "  -1.04%  "

$ws.Range("D39").Formula = "'329.34"
$ws.Range("E39").Value = "  +1.17%  "

$ws.Range("E40").Value = "  -1.00%  "

$ws.Range("D41").Formula = "'38.83"
$ws.Range("E41").Value = "  -2.46%  "

$ws.Range("D42").Formula = "'144.96"
$ws.Range("E42").Value = "  +2.98%  "

$ws.Range("E43").Value = "  -1.74%  "

$ws.Range("D44").Formula = "'20.09"
$ws.Range("E44").Value = "  +2.27%  "

$ws.Range("D45").Formula = "'0.0965"
$ws.Range("E45").Value = "  +0.35%  "

$ws.Range("E46").Value = "  -2.14%  "

$ws.Range("D47").Formula = "'0.576"
$ws.Range("E47").Value = "  +0.55%  "

$ws.Range("D49").Formula = "'11.03"
$ws.Range("E49").Value = "  -0.16%  "

$ws.Range("E50").Value = "  -3.24%  "

$ws.Range("B51").Value = "ZEEBU"
$ws.Range("C51").Value = "https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu"
$ws.Range("D51").Formula = "'4.65"
$ws.Range("E51").Value = "  -1.15%  "
